# Fruta / hortaliza, semanal
# Insert a new weekly record row for "Vega Modelo de Temuco" Kiwi data.
# This shifts existing rows 364-379 down to 365-380 and populates the
# newly created row 364 with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before current row 364 (shifts 364..379 -> 365..380)
$ws.Rows.Item(364).Insert()

# Populate the new row 364 with the new weekly data point
$ws.Cells.Item(364, 1).Value2 = 10
$ws.Cells.Item(364, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(364, 3).Value2 = "La Araucanía"
$ws.Cells.Item(364, 4).Value2 = 44568
$ws.Cells.Item(364, 5).Value2 = 9
$ws.Cells.Item(364, 6).Value2 = "Fruta"
$ws.Cells.Item(364, 7).Value2 = 100101
$ws.Cells.Item(364, 8).Value2 = "Berries"
$ws.Cells.Item(364, 9).Value2 = 100101007
$ws.Cells.Item(364, 10).Value2 = "Kiwi"
$ws.Cells.Item(364, 11).Value2 = "Hayward"
$ws.Cells.Item(364, 12).Value2 = "Segunda"
$ws.Cells.Item(364, 13).Value2 = 45
$ws.Cells.Item(364, 14).Value2 = 14000
$ws.Cells.Item(364, 15).Value2 = 14000
$ws.Cells.Item(364, 16).Value2 = 14000
$ws.Cells.Item(364, 17).Value2 = "`$/bandeja 18 kilos"
$ws.Cells.Item(364, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(364, 19).Value2 = 778
$ws.Cells.Item(364, 20).Value2 = 18
